$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# Row 8: swap the panel identifier text
$ws.Range("A8").Value = "FC32-1"

# Row 4: Update constant value and clear its highlight style (was style with border only -> default/no style)
$ws.Range("B4").Value = "NGC-1928/T959 OR TC-71688"
$ws.Range("B4").Style = "Normal"

# Row 8: CPU type values are swapped between columns C and F
$ws.Range("C8").Value = "CPU 801"
$ws.Range("F8").Value = "CPU 800"

# F8 needs to pick up the same formatting (fill/border/alignment) as the other
# data cells in that row (A8/C8), matching the style used for the rest of row 8.
$ws.Range("A8").Copy()
$ws.Range("F8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Updated measured current values
$ws.Range("G8").Value = 0.198
$ws.Range("H8").Value = 0.387

# Move the active selection to B4 to match where the edit was made
$ws.Range("B4").Select() | Out-Null
